# Updated symbol list on Thu Dec 15 22:24:19 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) — and two "Worstin24h"/"Bestin24h"
# suffix flags that moved rows in the "Volume(1h)" column (E) — to match
# the latest scrape of the crypto price table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells store plain numeric-looking strings (e.g. "258.32",
# "0.06080") as text, preserving trailing zeros / exact precision from the
# scrape. A naive "$ws.Range($cell).Value = $newValue" lets Excel's COM
# layer auto-coerce a numeric-looking string into a real number (dropping
# formatting such as trailing zeros) and also silently reassigns a new
# cell style. To keep both the literal text AND the original (default)
# style untouched, temporarily force the cell to Text format while we
# write, then restore its original style afterwards.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = $origStyle
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "258.00"
Set-TextValue "D3"  "22.75"
Set-TextValue "D4"  "6.166"
Set-TextValue "D5"  "0.06078"
Set-TextValue "D6"  "6.713"
Set-TextValue "D7"  "3.460"
Set-TextValue "D9"  "0.7966"
Set-TextValue "D10" "0.1581"
Set-TextValue "D11" "0.08049"
Set-TextValue "D12" "0.03345"
Set-TextValue "D14" "0.09298"
Set-TextValue "D15" "3.909"
Set-TextValue "D16" "0.001708"
Set-TextValue "D17" "0.04830"
Set-TextValue "D18" "0.0006147"
Set-TextValue "D19" "0.006245"
Set-TextValue "D20" "0.001102"
Set-TextValue "D21" "0.003383"
Set-TextValue "D23" "3.690"
Set-TextValue "D24" "2.261"
Set-TextValue "D25" "0.3356"
Set-TextValue "D27" "0.0003019"
Set-TextValue "D40" "0.04575"
Set-TextValue "D41" "0.007173"
Set-TextValue "D43" "0.1114"
Set-TextValue "D44" "0.009908"
Set-TextValue "D45" "0.002974"
Set-TextValue "D46" "0.00005990"
Set-TextValue "D48" "0.7509"
Set-TextValue "D49" "0.1085"

# --- Volume(1h) (column E) label updates ---
# Row 18 (One/ONE) picked up the "Worstin24h" suffix...
$ws.Range("E18").Value = "17OneONEWorstin24h"
# ...while row 21 (HotbitToken/HTB) lost it.
$ws.Range("E21").Value = "20HotbitTokenHTB"
